$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the organization website value in B10 from "www.stat.kg" to "www.stat.gov.kg"
$ws.Range("B10").Value = "www.stat.gov.kg"

# Move the active selection to B10 (matches the saved cursor position in the file)
$ws.Range("B10").Select() | Out-Null
